$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.365.34'
$ws.Range('E2').Value = '  +0.71%  '
$ws.Range('D3').Value = '2.761.70'
$ws.Range('E3').Value = '  +0.41%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '576.96'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '160.49'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.81%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.603'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.43%  '
$ws.Range('E9').Value = '  -1.44%  '
$ws.Range('E10').Value = '  +4.79%  '
$ws.Range('E11').Value = '  +2.70%  '
$ws.Range('E12').Value = '  -1.32%  '
$ws.Range('D13').Value = '3.251.29'
$ws.Range('E13').Value = '  +0.47%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.40'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.32%  '
$ws.Range('D15').Value = '63.995.28'
$ws.Range('E15').Value = '  +0.34%  '
$ws.Range('E16').Value = '  -2.10%  '
$ws.Range('D17').Value = '2.767.72'
$ws.Range('E17').Value = '  +0.54%  '
$ws.Range('E18').Value = '  -0.92%  '
$ws.Range('E19').Value = '  -2.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '358.42'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.90%  '
$ws.Range('E21').Value = '  -3.46%  '
$ws.Range('E22').Value = '  +0.74%  '
$ws.Range('E23').Value = '  -6.83%  '
$ws.Range('E24').Value = '  -1.73%  '
$ws.Range('E25').Value = '  -1.19%  '
$ws.Range('E26').Value = '  -0.80%  '
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('E28').Value = '  -1.34%  '
$ws.Range('B29').Value = 'Fetch.AI'
$ws.Range('C29').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.38'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +10.15%  '
$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.38'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.29%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.98'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.30%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '168.55'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.06%  '
$ws.Range('E33').Value = '  +2.98%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.99'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.73%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '20.22'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.79%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.84'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.82%  '
$ws.Range('E38').Value = '  -0.67%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '350.72'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.19%  '
$ws.Range('E40').Value = '  +3.69%  '
$ws.Range('E41').Value = '  -0.41%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '38.98'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.88%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '22.55'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.35%  '
$ws.Range('E44').Value = '  -2.18%  '
$ws.Range('E45').Value = '  -0.80%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '136.97'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.18%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.631'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.64%  '
$ws.Range('E48').Value = '  -2.13%  '
$ws.Range('E49').Value = '  -1.14%  '
$ws.Range('E50').Value = '  -0.15%  '
$ws.Range('D51').Value = '2.149.30'
$ws.Range('E51').Value = '  +1.33%  '
